$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric cells
$ws.Range("A2").Value = 182234
$ws.Range("B2").Value = 104451
$ws.Range("S2").Value = 5

# Text cells that must stay text even though they look numeric / date-like
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "200"

$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2010-08-26"

$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2010-08-26"

# Plain text cells
$ws.Range("C2").Value = "Godkänd baserat på observatörens uppgifter"
$ws.Range("J2").Value = "plantor/tuvor"
$ws.Range("P2").Value = "Mästocka skjutfält, Ringvägen (12), Hl"
$ws.Range("X2").Value = "N-Lah-0790"
$ws.Range("AW2").Value = "Halland Floraväktarna"
$ws.Range("AX2").Value = "Lars-Erik Magnusson"
$ws.Range("AY2").Value = "Floraväkteri Sverige"

# Cell removed entirely in the diff (biotope description no longer set)
$ws.Range("AI2").ClearContents()
